$a = "hello"
$b = "world"
Write-Host ($a + " " + $b)

$x = 1.1
Write-Host ("val: " + $x)
